$d = $word.ActiveDocument

# Mapping of old three-digit multiplication strings to new ones.
# Using literal (non-wildcard) Find/Execute for each unique cell value.
$replacements = @(
    @{ Old = "829×3=2487"; New = "316×3=948" },
    @{ Old = "270×9=2430"; New = "498×8=3984" },
    @{ Old = "700×7=4900"; New = "119×8=952" },
    @{ Old = "196×4=784"; New = "268×7=1876" },
    @{ Old = "264×2=528"; New = "751×5=3755" },
    @{ Old = "993×3=2979"; New = "501×6=3006" },
    @{ Old = "130×6=780"; New = "192×8=1536" },
    @{ Old = "326×8=2608"; New = "418×8=3344" },
    @{ Old = "251×7=1757"; New = "370×5=1850" },
    @{ Old = "273×4=1092"; New = "796×7=5572" },
    @{ Old = "900×7=6300"; New = "746×9=6714" },
    @{ Old = "284×7=1988"; New = "348×5=1740" },
    @{ Old = "677×4=2708"; New = "975×4=3900" },
    @{ Old = "748×3=2244"; New = "638×9=5742" },
    @{ Old = "850×7=5950"; New = "694×2=1388" },
    @{ Old = "424×2=848"; New = "462×3=1386" },
    @{ Old = "782×5=3910"; New = "923×5=4615" },
    @{ Old = "155×3=465"; New = "186×2=372" },
    @{ Old = "120×7=840"; New = "988×2=1976" },
    @{ Old = "396×3=1188"; New = "482×3=1446" },
    @{ Old = "911×7=6377"; New = "782×3=2346" },
    @{ Old = "140×2=280"; New = "349×9=3141" },
    @{ Old = "861×5=4305"; New = "502×8=4016" },
    @{ Old = "268×8=2144"; New = "439×9=3951" },
    @{ Old = "325×6=1950"; New = "894×6=5364" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
